# ---------------------------------------------------------------------------
# Applies the APP-000020 assets_liabilities.xlsx update:
#   - Summary sheet: new client name + refreshed income/assets/liabilities/
#     net-worth/ratio figures.
#   - Assets sheet: a new "Vehicles / Premium Car" row is inserted above the
#     existing "Liquid Assets" row, the savings balance drops to 5000, and
#     TOTAL ASSETS is recomputed.
#   - Liabilities sheet: a new "Auto Loans / Vehicle Loan 1" row is inserted
#     above the existing "Credit Cards" row, the credit-card figures change,
#     and TOTAL LIABILITIES is recomputed.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ===========================================================================
# 1) Summary sheet - simple in-place value updates (no structural changes)
# ===========================================================================
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("B3").Value = "Ahmed Al Mazrouei"   # Name
$summary.Range("B4").Value = 4751.78               # Monthly Income (AED)
$summary.Range("B6").Value = 222797                # Total Assets (AED)
$summary.Range("B7").Value = 175457                # Total Liabilities (AED)
$summary.Range("B8").Value = 47340                 # Net Worth (AED)
$summary.Range("B9").Value = 1.27                  # Asset/Liability Ratio

# ===========================================================================
# 2) Assets sheet - insert a new row above the existing data row
# ===========================================================================
$assets = $wb.Worksheets.Item("Assets")

# Shift "Liquid Assets" (row 2) and "TOTAL ASSETS" (row 3) down by inserting
# a fresh blank row at position 2; Excel carries each existing row's own
# formatting down with it, so row 3 keeps the "data row" style and row 4
# keeps the "TOTAL ASSETS" style automatically.
$assets.Rows.Item(2).Insert()

# The freshly inserted row 2 defaults to "copy formatting from the row
# above" (the header). Clear that first, then rebuild the plain data-row
# look (thin border, default font, "#,##0" value format) that the other
# data rows on this sheet use.
$newRow = $assets.Range("A2:C2")
$newRow.ClearFormats()
$assets.Range("A2:B2").Borders.LineStyle = 1
$assets.Range("C2").Borders.LineStyle = 1
$assets.Range("C2").NumberFormat = "#,##0"

$assets.Range("A2").Value = "Vehicles"
$assets.Range("B2").Value = "Premium Car"
$assets.Range("C2").Value = 217797

# Update the (now shifted-down) existing rows' figures.
$assets.Range("C3").Value = 5000      # Liquid Assets / Savings Account
$assets.Range("C4").Value = 222797    # TOTAL ASSETS

# ===========================================================================
# 3) Liabilities sheet - insert a new row above the existing data row
# ===========================================================================
$liab = $wb.Worksheets.Item("Liabilities")

$liab.Rows.Item(2).Insert()

$newRow2 = $liab.Range("A2:E2")
$newRow2.ClearFormats()
$liab.Range("A2:B2").Borders.LineStyle = 1
$liab.Range("E2").Borders.LineStyle = 1
$liab.Range("C2:D2").Borders.LineStyle = 1
$liab.Range("C2:D2").NumberFormat = "#,##0"

$liab.Range("A2").Value = "Auto Loans"
$liab.Range("B2").Value = "Vehicle Loan 1"
$liab.Range("C2").Value = 130678
$liab.Range("D2").Value = 3630
$liab.Range("E2").Value = 3

# Update the (now shifted-down) existing "Credit Cards" row and totals.
$liab.Range("C3").Value = 44779
$liab.Range("D3").Value = 2239
$liab.Range("C4").Value = 175457   # TOTAL LIABILITIES
